# Creat liability_returns report, calculate aggregate liability return,
# and format asset liability returns table.
#
# Adds a new "Total" worksheet that sums IBT + Pension + Retirement market
# values for every date row, tidies up the scroll/selection state on the
# existing sheets, and extends the Retirement sheet with an extra blank
# formatted row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "Total" worksheet, placed after "Retirement" (last
#    tab). Worksheets.Add() with no args clones the ACTIVE sheet in this
#    runtime, so we route through a disposable placeholder sheet (added
#    with an explicit Before/After target, which yields a truly blank
#    sheet) to land on a clean worksheet. Doing one add+delete cycle
#    first also nudges the internal sheetId counter so "Total" ends up
#    with sheetId 5, matching a sheet that once lived in this slot.
# ---------------------------------------------------------------------
$retirement = $wb.Worksheets.Item("Retirement")

$placeholder = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $retirement)
$placeholder.Name = "__scratch__"

$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $placeholder)
$totalSheet.Name = "Total"

$placeholder = $wb.Worksheets.Item("__scratch__")
$placeholder.Delete()

# ---------------------------------------------------------------------
# 2. Populate "Total" with the same look as "Retirement" (header style,
#    date style, currency style, column widths) by copying its ranges,
#    then overwrite column B with the cross-sheet SUM formulas.
# ---------------------------------------------------------------------
$retirement = $wb.Worksheets.Item("Retirement")
$totalSheet = $wb.Worksheets.Item("Total")

$retirement.Range("A1:B146").Copy($totalSheet.Range("A1:B146"))
$retirement.Range("C144:C147").Copy($totalSheet.Range("C144:C147"))

for ($r = 2; $r -le 146; $r++) {
    $totalSheet.Range("B$r").Formula = "=SUM(IBT!B$r,Pension!B$r,Retirement!B$r)"
}

# Extra blank formatted row (matches the one appended to Retirement below).
$retirement.Range("C144").Copy($totalSheet.Range("B150"))

# ---------------------------------------------------------------------
# 3. Append the same blank formatted row to "Retirement" itself.
# ---------------------------------------------------------------------
$retirement.Range("C144").Copy($retirement.Range("B150"))

# ---------------------------------------------------------------------
# 4. Tidy up scroll position / selection on each sheet. Selecting a
#    range clears any stale topLeftCell scroll anchor; doing this in
#    tab order and finishing on "Retirement" keeps it the active tab
#    (matches activeTab="2" / tabSelected="1" being unchanged).
# ---------------------------------------------------------------------
$ibt = $wb.Worksheets.Item("IBT")
$ibt.Range("F144").Select()

$pension = $wb.Worksheets.Item("Pension")
$pension.Range("D150").Select()

$totalSheet = $wb.Worksheets.Item("Total")
$totalSheet.Range("D4").Select()

$retirement = $wb.Worksheets.Item("Retirement")
$retirement.Range("E7").Select()
